# Apply updated retention metrics to Sheet1.
# Columns: A=cohort_year, B=period_index, C=num_customers, D=cohort_size, E=retention_rate

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 2020 cohort (rows 2-7), cohort_size stays 107 ---
$ws.Range("C3").Value = 65
$ws.Range("E3").Value = 0.6074766355140186

$ws.Range("C4").Value = 27
$ws.Range("E4").Value = 0.2523364485981308

$ws.Range("C5").Value = 18
$ws.Range("E5").Value = 0.1682242990654206

$ws.Range("C6").Value = 16
$ws.Range("E6").Value = 0.1495327102803738

$ws.Range("C7").Value = 15
$ws.Range("E7").Value = 0.1401869158878505

# --- 2021 cohort (rows 8-12), cohort_size grows from 417 to 419 ---
$ws.Range("C8").Value = 419
$ws.Range("D8").Value = 419

$ws.Range("C9").Value = 142
$ws.Range("D9").Value = 419
$ws.Range("E9").Value = 0.3389021479713604

$ws.Range("C10").Value = 106
$ws.Range("D10").Value = 419
$ws.Range("E10").Value = 0.2529832935560859

$ws.Range("C11").Value = 87
$ws.Range("D11").Value = 419
$ws.Range("E11").Value = 0.20763723150358

$ws.Range("C12").Value = 71
$ws.Range("D12").Value = 419
$ws.Range("E12").Value = 0.1694510739856802

# --- 2022 cohort (rows 13-16), cohort_size stays 193 ---
$ws.Range("C15").Value = 63
$ws.Range("E15").Value = 0.3264248704663212

$ws.Range("C16").Value = 59
$ws.Range("E16").Value = 0.3056994818652849

# --- 2023 cohort (rows 17-19), cohort_size shrinks from 124 to 123 ---
$ws.Range("C17").Value = 123
$ws.Range("D17").Value = 123

$ws.Range("C18").Value = 85
$ws.Range("D18").Value = 123
$ws.Range("E18").Value = 0.6910569105691057

$ws.Range("C19").Value = 65
$ws.Range("D19").Value = 123
$ws.Range("E19").Value = 0.5284552845528455

# --- 2024 cohort (rows 20-21), cohort_size grows from 204 to 206 ---
$ws.Range("C20").Value = 206
$ws.Range("D20").Value = 206

$ws.Range("C21").Value = 114
$ws.Range("D21").Value = 206
$ws.Range("E21").Value = 0.5533980582524272

# --- 2025 cohort (row 22), cohort_size grows from 61 to 62 ---
$ws.Range("C22").Value = 62
$ws.Range("D22").Value = 62
